# Daily attendance processing - 2025-10-18 19:15:44
# Normalize the "Recorded By" (column G) cell values: for any cell whose
# value is a comma-separated list of recorder names/emails, rotate the
# list so the last entry moves to the front (e.g. "A, B" -> "B, A").
# Single-value cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $rotated = @($parts[-1]) + $parts[0..($parts.Count - 2)]
            $newVal = [string]::Join(", ", $rotated)
            $cell.Value2 = $newVal
        }
    }
}
